$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rebuild the account-statement detail table (rows 16-42).
#    The source system re-sorted the rows: instead of grouping by worker
#    (Jonathan, Alfonso, Leonardo, Maickol) with each worker's own period
#    list, the new layout groups by period (2102..2112, 2201) and lists the
#    relevant workers under each period. A new worker (Maickol Andres Polo
#    Silva, CC 1002257296) now also has entries for periods 2102-2106.
# ---------------------------------------------------------------------------
$data = @(
    ,@(16, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2102', 36341)
    ,@(17, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2103', 36341)
    ,@(18, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2104', 36341)
    ,@(19, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2105', 36341)
    ,@(20, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2106', 36341)
    ,@(21, '1128055824', 'JONATHAN TERAN TORRES', '2107', 36341)
    ,@(22, '92500773', 'ALFONSO RAFAEL TERAN MONTES', '2107', 36341)
    ,@(23, '1002322447', 'LEONARDO ENRIQUE MARIMON SANDOVAL', '2107', 36341)
    ,@(24, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2107', 36341)
    ,@(25, '1128055824', 'JONATHAN TERAN TORRES', '2108', 36341)
    ,@(26, '92500773', 'ALFONSO RAFAEL TERAN MONTES', '2108', 36341)
    ,@(27, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2108', 36341)
    ,@(28, '1128055824', 'JONATHAN TERAN TORRES', '2109', 36341)
    ,@(29, '92500773', 'ALFONSO RAFAEL TERAN MONTES', '2109', 36341)
    ,@(30, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2109', 36341)
    ,@(31, '1128055824', 'JONATHAN TERAN TORRES', '2110', 36341)
    ,@(32, '92500773', 'ALFONSO RAFAEL TERAN MONTES', '2110', 36341)
    ,@(33, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2110', 36341)
    ,@(34, '1128055824', 'JONATHAN TERAN TORRES', '2111', 36341)
    ,@(35, '92500773', 'ALFONSO RAFAEL TERAN MONTES', '2111', 36341)
    ,@(36, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2111', 36341)
    ,@(37, '1128055824', 'JONATHAN TERAN TORRES', '2112', 36341)
    ,@(38, '92500773', 'ALFONSO RAFAEL TERAN MONTES', '2112', 36341)
    ,@(39, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2112', 36341)
    ,@(40, '1128055824', 'JONATHAN TERAN TORRES', '2201', 30284)
    ,@(41, '92500773', 'ALFONSO RAFAEL TERAN MONTES', '2201', 30284)
    ,@(42, '1002257296', 'MAICKOL ANDRES POLO SILVA', '2201', 30284)
)


foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value2 = "CC"
    $ws.Range("C$r").Value2 = $row[1]
    $ws.Range("D$r").Value2 = $row[2]
    $ws.Range("E$r").Value2 = $row[3]
    $ws.Range("F$r").Value2 = $row[4]
}

# ---------------------------------------------------------------------------
# 2) Header-block border touch-up: the vertical divider between the
#    logo/title cells (and between each label/value pair) was shifted to
#    sit on the true right-hand edge of the merged label cell instead of on
#    an interior cell, and the stray outer-right edge on the value cells
#    was dropped to match.
# ---------------------------------------------------------------------------
$loseRightBorder = @("B2","B3","B4","B5","B7","C7","B9","C9","B11","C11",
    "E7","F7","G7","H7","I7","J7","E9","F9","G9","H9","I9","J9",
    "E11","F11","G11","H11","I11","J11",
    "B13","C13","D13","F13","G13","H13","I13","J13")
foreach ($addr in $loseRightBorder) {
    $ws.Range($addr).Borders.Item(10).LineStyle = -4142
}

$gainRightBorder = @("C2","C3","C4","C5","D7","D9","D11","E13")
foreach ($addr in $gainRightBorder) {
    $ws.Range($addr).Borders.Item(10).LineStyle = 1
    $ws.Range($addr).Borders.Item(10).Weight = 2
}
